$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column H ("Save"), copying the formatting (style) from G1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Values for H2:H12 ("Save" flag per row)
$values = @(1, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
